$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.101.10'
$ws.Range("E2").Value = '  -1.80%  '
$ws.Range("D3").Value = '1.789.50'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '222.19'
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.29'
$ws.Range("E8").Value = '  -0.73%  '
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("E10").Value = '  -1.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0928'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").Value = '2.045.85'
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("D13").Value = '1.807.06'
$ws.Range("E13").Value = '  +0.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.95'
$ws.Range("E14").Value = '  -2.02%  '
$ws.Range("E15").Value = '  -2.08%  '
$ws.Range("D16").Value = '34.103.78'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.17'
$ws.Range("E17").Value = '  -3.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.07'
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.02'
$ws.Range("E19").Value = '  -4.30%  '
$ws.Range("E20").Value = '  -3.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.78'
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("E23").Value = '  -3.75%  '
$ws.Range("E24").Value = '  -1.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.15'
$ws.Range("E25").Value = '  -1.05%  '
$ws.Range("E26").Value = '  -0.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.06'
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("E28").Value = '  -2.05%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0520'
$ws.Range("E30").Value = '  -3.15%  '
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.67'
$ws.Range("E32").Value = '  -3.78%  '
$ws.Range("E33").Value = '  -3.57%  '
$ws.Range("E34").Value = '  -4.09%  '
$ws.Range("D35").Value = '1.396.96'
$ws.Range("E35").Value = '  -3.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.653'
$ws.Range("E36").Value = '  +1.95%  '
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("E38").Value = '  -3.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '79.67'
$ws.Range("E39").Value = '  -6.79%  '
$ws.Range("E40").Value = '  +0.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.921'
$ws.Range("E41").Value = '  -2.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.71'
$ws.Range("E42").Value = '  -3.04%  '
$ws.Range("E43").Value = '  +2.01%  '
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0496'
$ws.Range("E44").Value = '  +0.36%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.91'
$ws.Range("E45").Value = '  -1.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '107.51'
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("E47").Value = '  -1.03%  '
$ws.Range("D48").Value = '1.946.94'
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.05'
$ws.Range("E49").Value = '  -0.16%  '
$ws.Range("E50").Value = '  -0.09%  '
